# Update cryptos price/volume data per upstream refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" "43.107.80"
Set-TextCell "E2" "  +0.26%  "
Set-TextCell "D3" "2.328.83"
Set-TextCell "E3" "  +1.04%  "
Set-TextCell "E4" "  +0.04%  "
Set-TextCell "D5" "304.77"
Set-TextCell "E5" "  +1.66%  "
Set-TextCell "D6" "97.81"
Set-TextCell "E6" "  +0.39%  "
Set-TextCell "E7" "  -1.32%  "
Set-TextCell "D9" "0.503"
Set-TextCell "E9" "  -0.50%  "
Set-TextCell "D10" "35.58"
Set-TextCell "E10" "  -0.45%  "
Set-TextCell "D11" "19.41"
Set-TextCell "E11" "  +7.61%  "
Set-TextCell "D12" "0.0799"
Set-TextCell "E12" "  +1.40%  "
Set-TextCell "E13" "  +1.15%  "
Set-TextCell "E14" "  +1.76%  "
Set-TextCell "D15" "2.692.16"
Set-TextCell "E15" "  +1.13%  "
Set-TextCell "D16" "2.335.29"
Set-TextCell "E16" "  +1.38%  "
Set-TextCell "E17" "  +1.05%  "
Set-TextCell "D18" "43.024.65"
Set-TextCell "E18" "  +0.25%  "
Set-TextCell "D19" "12.53"
Set-TextCell "E19" "  -1.23%  "
Set-TextCell "D20" "0.0₃0902"
Set-TextCell "E20" "  -0.14%  "
Set-TextCell "E21" "  +0.69%  "
Set-TextCell "D22" "67.97"
Set-TextCell "E22" "  -0.06%  "
Set-TextCell "D23" "237.64"
Set-TextCell "E23" "  -1.06%  "
Set-TextCell "E24" "  +3.81%  "
Set-TextCell "E25" "  +0.06%  "
Set-TextCell "E26" "  +0.30%  "
Set-TextCell "D27" "24.93"
Set-TextCell "E27" "  -2.34%  "
Set-TextCell "D28" "166.42"
Set-TextCell "E28" "  +0.55%  "
Set-TextCell "E29" "  +1.93%  "
Set-TextCell "D30" "9.13"
Set-TextCell "E30" "  +0.89%  "
Set-TextCell "D31" "33.15"
Set-TextCell "E31" "  +0.16%  "
Set-TextCell "D33" "17.98"
Set-TextCell "E33" "  +5.79%  "
Set-TextCell "E35" "  -8.47%  "
Set-TextCell "B36" "WEMIXToken"
Set-TextCell "C36" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell "D36" "2.35"
Set-TextCell "E36" "  -1.33%  "
Set-TextCell "B37" "Hedera"
Set-TextCell "C37" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell "D37" "0.0696"
Set-TextCell "E37" "  +1.22%  "
Set-TextCell "E38" "  +0.27%  "
Set-TextCell "D39" "2.81"
Set-TextCell "E39" "  +2.68%  "
Set-TextCell "D40" "1.76"
Set-TextCell "E40" "  +0.20%  "
Set-TextCell "D41" "0.109"
Set-TextCell "E41" "  -0.42%  "
Set-TextCell "D42" "1.997.35"
Set-TextCell "E42" "  -0.64%  "
Set-TextCell "D43" "10.76"
Set-TextCell "E43" "  +6.15%  "
Set-TextCell "E44" "  -0.07%  "
Set-TextCell "D45" "18.08"
Set-TextCell "E45" "  +4.91%  "
Set-TextCell "D46" "2.09"
Set-TextCell "E46" "  -2.19%  "
Set-TextCell "D47" "2.78"
Set-TextCell "E47" "  -0.53%  "
Set-TextCell "D48" "2.559.36"
Set-TextCell "E48" "  +1.16%  "
Set-TextCell "E49" "  -0.04%  "
Set-TextCell "D50" "53.81"
Set-TextCell "E50" "  +0.40%  "
Set-TextCell "D51" "71.94"
Set-TextCell "E51" "  -0.51%  "
